$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.892.15"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "1.739.50"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'230.87"
$ws.Range("E5").Value = "  -2.94%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'0.5257"
$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").Value = "'0.2768"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "'39.49"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("D10").Value = "'0.06151"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").Value = "1.736.98"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").Value = "'0.07111"
$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").Value = "'15.18"
$ws.Range("E13").Value = "  -3.47%  "

$ws.Range("D14").Value = "'0.6451"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "'4.529"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "'77.03"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").Value = "'0.9992"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").Value = "25.868.45"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").Value = "'0.000006687"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").Value = "1.959.77"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").Value = "'4.276"
$ws.Range("E23").Value = "  +4.94%  "

$ws.Range("D24").Value = "'8.803"
$ws.Range("E24").Value = "  +4.09%  "

$ws.Range("D25").Value = "'5.168"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").Value = "'140.41"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("D28").Value = "'15.19"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'1.794"
$ws.Range("E29").Value = "  -3.10%  "

$ws.Range("D30").Value = "'102.35"
$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").Value = "'3.739"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "'3.592"
$ws.Range("E33").Value = "  +4.02%  "

$ws.Range("D34").Value = "'0.04509"
$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("D35").Value = "'2.609"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9738"
$ws.Range("E36").Value = "  -2.94%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6203"
$ws.Range("E37").Value = "  +2.52%  "

$ws.Range("E38").Value = "  -2.72%  "

$ws.Range("D39").Value = "'0.01583"
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.910"
$ws.Range("E40").Value = "  -3.68%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'0.9992"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").Value = "'99.97"
$ws.Range("E42").Value = "  -2.67%  "

$ws.Range("D43").Value = "'0.3861"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").Value = "'0.7291"
$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("D45").Value = "'5.014"
$ws.Range("E45").Value = "  +1.69%  "

$ws.Range("D46").Value = "'0.05319"
$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").Value = "'6.220"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("D49").Value = "'53.58"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "'30.12"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").Value = "'7.653"
$ws.Range("E51").Value = "  +3.13%  "
